# Robot Facturas - recopilación de datos y asignación a celdas Excel
# Populates the "Seguimiento_Facturas" sheet with the invoice tracking table:
# headers (Número de Factura, Fecha, Importe) + 10 rows of invoice data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1,1).Value = "Número de Factura"
$ws.Cells.Item(1,2).Value = "Fecha"
$ws.Cells.Item(1,3).Value = "Importe"

# The invoice-number column holds values that look numeric ("1", "2", ... "10");
# the source data keeps them as text, so the cells must stay as text instead of
# being auto-converted to numbers by Excel.
$ws.Range("A2:A11").NumberFormat = "@"

$invoiceNumbers = @("1","10","2","3","4","5","6","7","8","9")
$fechas = @("31/06/2021","31/06/2021","31/06/2021","31/06/2021","31/06/2021","31/06/2021","31/06/2021","31/06/2021","31/06/2021","31/06/2021")
$importes = @(
    "Total impuestos 54,60",
    "Total impuestos 42,00",
    "Total impuestos 29,40",
    "Total impuestos 63,00",
    "Total impuestos 71,40",
    "Total impuestos 42,00",
    "Total impuestos 67,20",
    "Total impuestos 46,20",
    "Total impuestos 128,10",
    "Total impuestos 54,60"
)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $invoiceNumbers[$i]
    $ws.Cells.Item($row, 2).Value = $fechas[$i]
    $ws.Cells.Item($row, 3).Value = $importes[$i]
}
